$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 124, pushing existing rows 124-150
# down to 125-151 (and extending the used range to A1:R151).
$ws.Rows.Item(124).Insert()

# Populate the newly-inserted row 124 with its data.
$ws.Range("A124").Value = 8
$ws.Range("B124").Value = "Terminal La Palmera de La Serena"
$ws.Range("C124").Value = "Coquimbo"
$ws.Range("D124").Value = 44943
$ws.Range("E124").Value = 4
$ws.Range("F124").Value = 100112052
$ws.Range("G124").Value = "Albahaca"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 400
$ws.Range("K124").Value = 4000
$ws.Range("L124").Value = 4500
$ws.Range("M124").Value = 4250
$ws.Range("N124").Value = "$/docena de matas"
$ws.Range("O124").Value = "Provincia del Elquí"
$ws.Range("P124").Value = 708
$ws.Range("Q124").Value = 6
$ws.Range("R124").Value = "Hortaliza"
